$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The worksheet is protected; temporarily unprotect so the cell values
# below (which live on locked cells) can be updated, then restore
# protection afterwards.
$ws.Unprotect()

# Refresh the "as of" date in the confidentiality / disclosure footer (A11)
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-04-26 for illustrative purposes only and are subject to change."
# Re-assigning multi-line text can make Excel auto-expand the row height;
# re-fit it back down so row 11 keeps its original (default) height.
$ws.Rows.Item(11).AutoFit()

# Refresh the Weight (D) and Percent Change (E) figures for each holding
$ws.Range("D2").Value = 0.4899514472662832
$ws.Range("E2").Value = -0.0004810004810003043

$ws.Range("D3").Value = 0.2509798730535962
$ws.Range("E3").Value = 0.004585184123800046

$ws.Range("D4").Value = 0.09980760247587951
$ws.Range("E4").Value = 0.006520164211543111

$ws.Range("D5").Value = 0.1021478205213058
$ws.Range("E5").Value = 0.002621722846442021

$ws.Range("D6").Value = 0.02943923141990444
$ws.Range("E6").Value = 0.001949317738791478

$ws.Range("D7").Value = 0.02767402526303063
$ws.Range("E7").Value = 0.007831508515815155

$ws.Range("E8").Value = 0.002107803060532154

# Restore sheet protection to its prior (protected) state.
$ws.Protect()
